$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (border style) from the last existing data row (125)
# down to the new rows (126:132) so the new rows match the look of the
# existing table before filling in values.
$ws.Range("A125:G125").Copy()
$ws.Range("A126:G132").PasteSpecial(-4122)

# Row 126 - Al rain
$ws.Range("A126").Value = "Al rain"
$ws.Range("B126").Value = "Al rain"
$ws.Range("C126").Value = "الرين"
$ws.Range("D126").Value = 23.542764
$ws.Range("E126").Value = 45.515282
$ws.Range("F126").Value = "منطقة الرياض"
$ws.Range("G126").Value = "وسط المملكة"

# Row 127 - Al Sulayyil
$ws.Range("A127").Value = "Al Sulayyil"
$ws.Range("B127").Value = "Al Sulayyil"
$ws.Range("C127").Value = "السليل"
$ws.Range("D127").Value = 20.46646
$ws.Range("E127").Value = 45.56256
$ws.Range("F127").Value = "منطقة الرياض"
$ws.Range("G127").Value = "وسط المملكة"

# Row 128 - Rowaidat Alard
$ws.Range("A128").Value = "Rowaidat Alard"
$ws.Range("B128").Value = "Rowaidat Alard"
$ws.Range("C128").Value = "رويضة العرض"
$ws.Range("D128").Value = 23.774209
$ws.Range("E128").Value = 44.763002
$ws.Range("F128").Value = "منطقة الرياض"
$ws.Range("G128").Value = "وسط المملكة"

# Row 129 - Al Gara
$ws.Range("A129").Value = "Al Gara"
$ws.Range("B129").Value = "Al Gara"
$ws.Range("C129").Value = "القرى"
$ws.Range("D129").Value = 18.241785
$ws.Range("E129").Value = 42.488838
$ws.Range("F129").Value = "منطقة الباحة"
$ws.Range("G129").Value = "جنوب المملكة"

# Row 130 - Taima
$ws.Range("A130").Value = "Taima"
$ws.Range("B130").Value = "Taima"
$ws.Range("C130").Value = "تيماء"
$ws.Range("D130").Value = 27.61013
$ws.Range("E130").Value = 38.521644
$ws.Range("F130").Value = "منطقة تبوك"
$ws.Range("G130").Value = "شمال المملكة"

# Row 131 - Wethylan
$ws.Range("A131").Value = "Wethylan"
$ws.Range("B131").Value = "Wethylan"
$ws.Range("C131").Value = "وثيلان"
$ws.Range("D131").Value = 24.540605
$ws.Range("E131").Value = 46.693131
$ws.Range("F131").Value = "منطقة الرياض"
$ws.Range("G131").Value = "وسط المملكة"

# Row 132 - Marat
$ws.Range("A132").Value = "Marat"
$ws.Range("B132").Value = "Marat"
$ws.Range("C132").Value = "مرات"
$ws.Range("D132").Value = 25.069855
$ws.Range("E132").Value = 45.464037
$ws.Range("F132").Value = "منطقة الرياض"
$ws.Range("G132").Value = "وسط المملكة"

# Update the selection shown for the sheet's pane to match the new extent.
$ws.Range("A1:G132").Select() | Out-Null
